# Rename the single worksheet from "Sheet1" to "Program Info".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Program Info"

# Move the active selection on that sheet from E10 to E22.
$ws.Range("E22").Select()
